# edit.ps1
# Applies the LOB1229.docx content reshuffle described by the commit diff.
#
# The edit relocates several paragraph bodies to different headings without
# touching the paragraph/run formatting that is already in place (styles, bold
# labels, italics, manual line breaks). Word exposes a manual line break
# (<w:br/>) as Chr(11) inside Range/Find text, so plain Find.Execute text
# search/replace (MatchWildcards=$false) is sufficient to relocate every block,
# including the ones whose text spans a <w:br/>.
#
# Several destination texts are sourced from OTHER paragraphs in the same
# document (the moves form closed rotation cycles), so doing a naive sequential
# find/replace would let a later step match text a previous step just wrote.
# To make the result independent of execution order we first swap every source
# block for a short unique placeholder token (phase 1), and only then expand
# every placeholder into its final text (phase 2).

$d = $word.ActiveDocument
$vt = [char]11  # manual line break (<w:br/>) marker in Word range/find text

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $find"
    }
}

# ---------- Phase 1: move every source block to a unique placeholder token ----------
Replace-Text ("O Trabalho de Graduação 2 tem por objetivo a integração, o aprofundamento e aplicação dos conhecimentos adquiridos ao longo do curso, preparando e desenvolvendo a capacidade do aluno para a realização de atividades, que fazem parte do perfil de atuação profissional do engenheiro ambiental. O desenvolvimento do Trabalho de Graduação 2 deverá seguir o modelo escolhido no Trabalho de Graduação 1, e permitirá o uso de tecnologias digitais ou outras metodologias para desenvolvimento de conteúdo ou produto aplicável que utilize conteúdo da Engenharia Ambiental ou correlatas.") "@@TOK_PA@@"
Replace-Text ("Graduation Work 2 aims to integrate, deepen and apply the knowledge acquired throughout the course, preparing and developing the student's ability to carry out activities, which are part of the professional performance profile of the environmental engineer. The development of Graduate Work 2 should follow the model chosen in Graduate Work 1, and will allow the use of digital technologies or other methodologies for the development of content or applicable product that uses Environmental Engineering or related content.") "@@TOK_PB@@"
Replace-Text ("4780627 - Ana Lucia Gabas Ferreira") "@@TOK_ANA@@"
Replace-Text ("O aluno deverá desenvolver o Trabalho de Graduação 2 seguindo o modelo escolhido e já utilizado no Trabalho de Graduação 1. Diante da escolha do modelo, o aluno deverá desenvolver conteúdo científico ou produto aplicável, ambos relacionados a Engenharia Ambiental, e ao final, o Trabalho de Graduação deverá ser avaliado por banca avaliadora.") "@@TOK_PC@@"
Replace-Text ("The student must develop the Graduate Work 2 following the model chosen and already used in the Graduate Work 1. When choosing the model, the student must develop scientific content or applicable product, both related to Environmental Engineering, and at the end, the Graduate Work must be evaluated by an evaluating panel.") "@@TOK_ITALIC_C@@"
Replace-Text ("Para ambos os modelos (artigo ou produto): O aluno deverá dar continuidade ao desenvolvimento do Trabalho iniciado no Trabalho de Graduação 1. O programa da disciplina é constituído pelas seguintes etapas: 1) Desenvolvimento do tema com base nas atividades de cronograma aprovados no Trabalho de Graduação 1. 2) Desenvolvimento do texto final, conforme modelo fornecido pelos responsáveis da disciplina. 3) Entrega da versão final do texto, com aprovação do orientador e com a indicação da banca de avaliação. 4) Avaliação e atribuição de notas pela banca avaliadora, que pode ou não ser a mesma composta anteriormente no Trabalho de Graduação 1, a critério do orientador.") "@@TOK_PD@@"
Replace-Text ("Para ambos os modelos (artigo ou produto): O aluno deverá apresentar um artigo científico seguindo modelo fornecido pelos responsáveis da disciplina. A versão final do texto deverá ser aprovada pelo orientador no ato da submissão e deverá ser apresentado para banca de avaliação composta por dois doutores ou especialistas na área do projeto. Alternativamente, pode-se apresentar o texto com formatação de outra revista desde que seja anexado o comprovante de submissão do artigo, nesse caso, o texto apresentado pode seguir as regras de formatação da revista escolhida. A revista escolhida deve ser indexada por algum sistema de base de dados (Web of Science, Scopus, Scielo, Cinahl, Medline, etc). Alternativamente, artigos aprovados em revistas da área de estudo, até a semana anterior a apresentação, desobrigam o aluno a apresentar o trabalho para a banca e nesse caso, a entrega do artigo deve ser acompanhada pela comprovação do aceite do trabalho. A revista escolhida deve ser indexada por algum sistema de base de dados (Web of Science, Scopus, Scielo, Cinahl, Medline, etc)" + $vt) "@@TOK_PE_SLOT@@"
Replace-Text ("Avaliação e emissão de parecer pela banca avaliadora e pelo orientador, com atribuição de nota única final." + $vt + "Fica sob responsabilidade do orientador a verificação de ocorrência de plágio utilizando software apropriado e avaliação em Comitê de Ética, quando exigido, via cadastro na Plataforma Brasil." + $vt) "@@TOK_PF_SLOT@@"
Replace-Text ("Não há.") "@@TOK_PG_SLOT@@"
Replace-Text ("A ser definido no decorrer de cada projeto") "@@TOK_PH@@"

# ---------- Phase 2: expand every placeholder into its final text ----------
Replace-Text "@@TOK_PA@@" ("O aluno deverá desenvolver o Trabalho de Graduação 2 seguindo o modelo escolhido e já utilizado no Trabalho de Graduação 1. Diante da escolha do modelo, o aluno deverá desenvolver conteúdo científico ou produto aplicável, ambos relacionados a Engenharia Ambiental, e ao final, o Trabalho de Graduação deverá ser avaliado por banca avaliadora.")
Replace-Text "@@TOK_PB@@" ("The student must develop the Graduate Work 2 following the model chosen and already used in the Graduate Work 1. When choosing the model, the student must develop scientific content or applicable product, both related to Environmental Engineering, and at the end, the Graduate Work must be evaluated by an evaluating panel.")
Replace-Text "@@TOK_ANA@@" ("O Trabalho de Graduação 2 tem por objetivo a integração, o aprofundamento e aplicação dos conhecimentos adquiridos ao longo do curso, preparando e desenvolvendo a capacidade do aluno para a realização de atividades, que fazem parte do perfil de atuação profissional do engenheiro ambiental. O desenvolvimento do Trabalho de Graduação 2 deverá seguir o modelo escolhido no Trabalho de Graduação 1, e permitirá o uso de tecnologias digitais ou outras metodologias para desenvolvimento de conteúdo ou produto aplicável que utilize conteúdo da Engenharia Ambiental ou correlatas.")
Replace-Text "@@TOK_PC@@" ("Para ambos os modelos (artigo ou produto): O aluno deverá dar continuidade ao desenvolvimento do Trabalho iniciado no Trabalho de Graduação 1. O programa da disciplina é constituído pelas seguintes etapas: 1) Desenvolvimento do tema com base nas atividades de cronograma aprovados no Trabalho de Graduação 1. 2) Desenvolvimento do texto final, conforme modelo fornecido pelos responsáveis da disciplina. 3) Entrega da versão final do texto, com aprovação do orientador e com a indicação da banca de avaliação. 4) Avaliação e atribuição de notas pela banca avaliadora, que pode ou não ser a mesma composta anteriormente no Trabalho de Graduação 1, a critério do orientador.")
Replace-Text "@@TOK_ITALIC_C@@" ("Graduation Work 2 aims to integrate, deepen and apply the knowledge acquired throughout the course, preparing and developing the student's ability to carry out activities, which are part of the professional performance profile of the environmental engineer. The development of Graduate Work 2 should follow the model chosen in Graduate Work 1, and will allow the use of digital technologies or other methodologies for the development of content or applicable product that uses Environmental Engineering or related content.")
Replace-Text "@@TOK_PD@@" ("Para ambos os modelos (artigo ou produto): O aluno deverá apresentar um artigo científico seguindo modelo fornecido pelos responsáveis da disciplina. A versão final do texto deverá ser aprovada pelo orientador no ato da submissão e deverá ser apresentado para banca de avaliação composta por dois doutores ou especialistas na área do projeto. Alternativamente, pode-se apresentar o texto com formatação de outra revista desde que seja anexado o comprovante de submissão do artigo, nesse caso, o texto apresentado pode seguir as regras de formatação da revista escolhida. A revista escolhida deve ser indexada por algum sistema de base de dados (Web of Science, Scopus, Scielo, Cinahl, Medline, etc). Alternativamente, artigos aprovados em revistas da área de estudo, até a semana anterior a apresentação, desobrigam o aluno a apresentar o trabalho para a banca e nesse caso, a entrega do artigo deve ser acompanhada pela comprovação do aceite do trabalho. A revista escolhida deve ser indexada por algum sistema de base de dados (Web of Science, Scopus, Scielo, Cinahl, Medline, etc)")
Replace-Text "@@TOK_PE_SLOT@@" ("Avaliação e emissão de parecer pela banca avaliadora e pelo orientador, com atribuição de nota única final." + $vt + "Fica sob responsabilidade do orientador a verificação de ocorrência de plágio utilizando software apropriado e avaliação em Comitê de Ética, quando exigido, via cadastro na Plataforma Brasil." + $vt)
Replace-Text "@@TOK_PF_SLOT@@" ("Não há." + $vt)
Replace-Text "@@TOK_PG_SLOT@@" ("A ser definido no decorrer de cada projeto")
Replace-Text "@@TOK_PH@@" ("4780627 - Ana Lucia Gabas Ferreira")

Write-Output "done"
